$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "60-57=3"
$t.Cell(1, 2).Range.Text = "7+55=62"
$t.Cell(1, 3).Range.Text = "48+17=65"
$t.Cell(1, 4).Range.Text = "59+16=75"
$t.Cell(1, 5).Range.Text = "84-16=68"
$t.Cell(2, 1).Range.Text = "64-19=45"
$t.Cell(2, 2).Range.Text = "93-44=49"
$t.Cell(2, 3).Range.Text = "63-14=49"
$t.Cell(2, 4).Range.Text = "7+77=84"
$t.Cell(2, 5).Range.Text = "58+36=94"
$t.Cell(3, 1).Range.Text = "9+39=48"
$t.Cell(3, 2).Range.Text = "45+46=91"
$t.Cell(3, 3).Range.Text = "95-67=28"
$t.Cell(3, 4).Range.Text = "23+19=42"
$t.Cell(3, 5).Range.Text = "12+79=91"
$t.Cell(4, 1).Range.Text = "67+17=84"
$t.Cell(4, 2).Range.Text = "74-26=48"
$t.Cell(4, 3).Range.Text = "49+42=91"
$t.Cell(4, 4).Range.Text = "87-48=39"
$t.Cell(4, 5).Range.Text = "71-13=58"
$t.Cell(5, 1).Range.Text = "50-38=12"
$t.Cell(5, 2).Range.Text = "71-42=29"
$t.Cell(5, 3).Range.Text = "7+18=25"
$t.Cell(5, 4).Range.Text = "75-27=48"
$t.Cell(5, 5).Range.Text = "9+87=96"
$t.Cell(6, 1).Range.Text = "18+58=76"
$t.Cell(6, 2).Range.Text = "55-18=37"
$t.Cell(6, 3).Range.Text = "9+35=44"
$t.Cell(6, 4).Range.Text = "81-9=72"
$t.Cell(6, 5).Range.Text = "81-3=78"
$t.Cell(7, 1).Range.Text = "90-89=1"
$t.Cell(7, 2).Range.Text = "95-46=49"
$t.Cell(7, 3).Range.Text = "18+18=36"
$t.Cell(7, 4).Range.Text = "42-39=3"
$t.Cell(7, 5).Range.Text = "18+64=82"
$t.Cell(8, 1).Range.Text = "47+44=91"
$t.Cell(8, 2).Range.Text = "13+19=32"
$t.Cell(8, 3).Range.Text = "82-78=4"
$t.Cell(8, 4).Range.Text = "16+36=52"
$t.Cell(8, 5).Range.Text = "71-18=53"
$t.Cell(9, 1).Range.Text = "15+39=54"
$t.Cell(9, 2).Range.Text = "76-57=19"
$t.Cell(9, 3).Range.Text = "19+63=82"
$t.Cell(9, 4).Range.Text = "19+22=41"
$t.Cell(9, 5).Range.Text = "9+16=25"
$t.Cell(10, 1).Range.Text = "6+59=65"
$t.Cell(10, 2).Range.Text = "90-16=74"
$t.Cell(10, 3).Range.Text = "15+48=63"
$t.Cell(10, 4).Range.Text = "19+66=85"
$t.Cell(10, 5).Range.Text = "61-9=52"
$t.Cell(11, 1).Range.Text = "61-49=12"
$t.Cell(11, 2).Range.Text = "80-7=73"
$t.Cell(11, 3).Range.Text = "79+3=82"
$t.Cell(11, 4).Range.Text = "83-48=35"
$t.Cell(11, 5).Range.Text = "84-76=8"
$t.Cell(12, 1).Range.Text = "18+13=31"
$t.Cell(12, 2).Range.Text = "39+47=86"
$t.Cell(12, 3).Range.Text = "58-9=49"
$t.Cell(12, 4).Range.Text = "51-9=42"
$t.Cell(12, 5).Range.Text = "28+8=36"
$t.Cell(13, 1).Range.Text = "90-29=61"
$t.Cell(13, 2).Range.Text = "96-49=47"
$t.Cell(13, 3).Range.Text = "63-29=34"
$t.Cell(13, 4).Range.Text = "49+5=54"
$t.Cell(13, 5).Range.Text = "69+9=78"
$t.Cell(14, 1).Range.Text = "21-4=17"
$t.Cell(14, 2).Range.Text = "38-19=19"
$t.Cell(14, 3).Range.Text = "69+27=96"
$t.Cell(14, 4).Range.Text = "93-88=5"
$t.Cell(14, 5).Range.Text = "18+66=84"
$t.Cell(15, 1).Range.Text = "73-38=35"
$t.Cell(15, 2).Range.Text = "48+6=54"
$t.Cell(15, 3).Range.Text = "57+6=63"
$t.Cell(15, 4).Range.Text = "49+22=71"
$t.Cell(15, 5).Range.Text = "15+57=72"
$t.Cell(16, 1).Range.Text = "41-7=34"
$t.Cell(16, 2).Range.Text = "82-8=74"
$t.Cell(16, 3).Range.Text = "96-79=17"
$t.Cell(16, 4).Range.Text = "31-26=5"
$t.Cell(16, 5).Range.Text = "25+7=32"
$t.Cell(17, 1).Range.Text = "12-6=6"
$t.Cell(17, 2).Range.Text = "58+37=95"
$t.Cell(17, 3).Range.Text = "8+88=96"
$t.Cell(17, 4).Range.Text = "73-26=47"
$t.Cell(17, 5).Range.Text = "19+75=94"
$t.Cell(18, 1).Range.Text = "93-85=8"
$t.Cell(18, 2).Range.Text = "44+37=81"
$t.Cell(18, 3).Range.Text = "9+53=62"
$t.Cell(18, 4).Range.Text = "29+44=73"
$t.Cell(18, 5).Range.Text = "19+42=61"
$t.Cell(19, 1).Range.Text = "44-16=28"
$t.Cell(19, 2).Range.Text = "85-29=56"
$t.Cell(19, 3).Range.Text = "26+29=55"
$t.Cell(19, 4).Range.Text = "72-35=37"
$t.Cell(19, 5).Range.Text = "60-16=44"
$t.Cell(20, 1).Range.Text = "39+37=76"
$t.Cell(20, 2).Range.Text = "34+17=51"
$t.Cell(20, 3).Range.Text = "35+47=82"
$t.Cell(20, 4).Range.Text = "67-39=28"
$t.Cell(20, 5).Range.Text = "39+7=46"
